# Apply the edit described by the diff:
# - Insert 3 new rows at row 859 (shifting the existing data block down by 3 rows,
#   from A1:R957 to A1:R960).
# - Populate the 3 newly inserted rows (859-861) with new sampling data for
#   date 45132 (2023-07-25): "Morrón rojo", "Zafiro rojo" and "Zafiro verde".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 859, pushing existing rows 859..957 down to 862..960.
$ws.Range("A859:A861").EntireRow.Insert()

# Row 859: Morrón rojo
$ws.Cells.Item(859, 1).Value2 = 5
$ws.Cells.Item(859, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(859, 3).Value2 = "Maule"
$ws.Cells.Item(859, 4).Value2 = 45132
$ws.Cells.Item(859, 5).Value2 = 7
$ws.Cells.Item(859, 6).Value2 = 100112002
$ws.Cells.Item(859, 7).Value2 = "Pimiento"
$ws.Cells.Item(859, 8).Value2 = "Morrón rojo"
$ws.Cells.Item(859, 9).Value2 = "Primera"
$ws.Cells.Item(859, 10).Value2 = 300
$ws.Cells.Item(859, 11).Value2 = 12000
$ws.Cells.Item(859, 12).Value2 = 12000
$ws.Cells.Item(859, 13).Value2 = 12000
$ws.Cells.Item(859, 14).Value2 = "`$/caja 20 kilos"
$ws.Cells.Item(859, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(859, 16).Value2 = 600
$ws.Cells.Item(859, 17).Value2 = 20
$ws.Cells.Item(859, 18).Value2 = "Hortaliza"

# Row 860: Zafiro rojo
$ws.Cells.Item(860, 1).Value2 = 5
$ws.Cells.Item(860, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(860, 3).Value2 = "Maule"
$ws.Cells.Item(860, 4).Value2 = 45132
$ws.Cells.Item(860, 5).Value2 = 7
$ws.Cells.Item(860, 6).Value2 = 100112002
$ws.Cells.Item(860, 7).Value2 = "Pimiento"
$ws.Cells.Item(860, 8).Value2 = "Zafiro rojo"
$ws.Cells.Item(860, 9).Value2 = "Primera"
$ws.Cells.Item(860, 10).Value2 = 200
$ws.Cells.Item(860, 11).Value2 = 18000
$ws.Cells.Item(860, 12).Value2 = 18000
$ws.Cells.Item(860, 13).Value2 = 18000
$ws.Cells.Item(860, 14).Value2 = "`$/caja 15 kilos"
$ws.Cells.Item(860, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(860, 16).Value2 = 1200
$ws.Cells.Item(860, 17).Value2 = 15
$ws.Cells.Item(860, 18).Value2 = "Hortaliza"

# Row 861: Zafiro verde
$ws.Cells.Item(861, 1).Value2 = 5
$ws.Cells.Item(861, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(861, 3).Value2 = "Maule"
$ws.Cells.Item(861, 4).Value2 = 45132
$ws.Cells.Item(861, 5).Value2 = 7
$ws.Cells.Item(861, 6).Value2 = 100112002
$ws.Cells.Item(861, 7).Value2 = "Pimiento"
$ws.Cells.Item(861, 8).Value2 = "Zafiro verde"
$ws.Cells.Item(861, 9).Value2 = "Primera"
$ws.Cells.Item(861, 10).Value2 = 300
$ws.Cells.Item(861, 11).Value2 = 11000
$ws.Cells.Item(861, 12).Value2 = 11000
$ws.Cells.Item(861, 13).Value2 = 11000
$ws.Cells.Item(861, 14).Value2 = "`$/caja 15 kilos"
$ws.Cells.Item(861, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(861, 16).Value2 = 733
$ws.Cells.Item(861, 17).Value2 = 15
$ws.Cells.Item(861, 18).Value2 = "Hortaliza"
